# NATMI LR-pairs output refresh: Ntn4-Dcc.xlsx was regenerated with updated
# TPM input, which changes the Ligand/Receptor/Edge expression + specificity
# statistics for all rows. Write the recomputed values directly into the
# cached cells (the source file stores literal numbers, not formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sender ECs -> Target ECs)
$ws.Range("G2").Value = 2.721294
$ws.Range("H2").Value = 8.163882000000001
$ws.Range("I2").Value = 0.03340571984979829
$ws.Range("J2").Value = 0.03340571984979829
$ws.Range("M2").Value = 0.092904
$ws.Range("Q2").Value = 0.252819097776
$ws.Range("R2").Value = 2.275371879984001
$ws.Range("S2").Value = 0.03197741109137893
$ws.Range("T2").Value = 0.03197741109137893

# Row 3 (Sender ECs -> Target MuSCs)
$ws.Range("G3").Value = 2.721294
$ws.Range("H3").Value = 8.163882000000001
$ws.Range("I3").Value = 0.03340571984979829
$ws.Range("J3").Value = 0.03340571984979829
$ws.Range("O3").Value = 0.04275641311851519
$ws.Range("P3").Value = 0.04275641311851518
$ws.Range("Q3").Value = 0.011292463002
$ws.Range("R3").Value = 0.101632167018
$ws.Range("S3").Value = 0.001428308758419359
$ws.Range("T3").Value = 0.001428308758419359

# Row 4 (Sender FAPs -> Target ECs)
$ws.Range("I4").Value = 0.427684027063558
$ws.Range("J4").Value = 0.427684027063558
$ws.Range("M4").Value = 0.092904
$ws.Range("Q4").Value = 3.236771736744
$ws.Range("R4").Value = 29.130945630696
$ws.Range("S4").Value = 0.4093977921182382
$ws.Range("T4").Value = 0.4093977921182382

# Row 5 (Sender FAPs -> Target MuSCs)
$ws.Range("I5").Value = 0.427684027063558
$ws.Range("J5").Value = 0.427684027063558
$ws.Range("O5").Value = 0.04275641311851519
$ws.Range("P5").Value = 0.04275641311851518
$ws.Range("S5").Value = 0.01828623494531972
$ws.Range("T5").Value = 0.01828623494531971

# Row 6 (Sender MuSCs -> Target ECs)
$ws.Range("I6").Value = 0.5389102530866438
$ws.Range("J6").Value = 0.5389102530866438
$ws.Range("M6").Value = 0.092904
$ws.Range("Q6").Value = 4.078547164384
$ws.Range("R6").Value = 36.706924479456
$ws.Range("S6").Value = 0.5158683836718676
$ws.Range("T6").Value = 0.5158683836718676

# Row 7 (Sender MuSCs -> Target MuSCs)
$ws.Range("I7").Value = 0.5389102530866438
$ws.Range("J7").Value = 0.5389102530866438
$ws.Range("O7").Value = 0.04275641311851519
$ws.Range("P7").Value = 0.04275641311851518
$ws.Range("S7").Value = 0.02304186941477612
$ws.Range("T7").Value = 0.02304186941477611
